# Apply the "Add files via upload" edit to 17-1-2.xlsx
#
# Net effect of the source diff:
#   * Cell B10 text is updated from
#       "MF KR: www.minfin.kg; \nNSC KR: www.stat.kg"
#     to
#       "MF KR: www.minfin.kg; \nNSC KR:www.stat.gov.kg"
#     (the old shared string is dropped and the new text is appended at the
#     end of the shared-strings table, which also re-slots B10 onto a
#     dedicated style/font entry in the real Excel session that produced the
#     diff).
#   * The worksheet's active selection moves from B2 to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B10")
$cell.Value = "MF KR: www.minfin.kg; `nNSC KR:www.stat.gov.kg"

# Give the edited cell its own font instance (mirrors the distinct font /
# cell-style slot that Excel allocated for this cell in the authored file).
$cell.Font.Name = "Calibri"

# Move the selection to match the saved workbook state (B9 selected).
$ws.Range("B9").Select()
